# Replace the single data row (row 2) with the new exposure-site record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Southern Cross"
$ws.Range("B2").Value = "Metro trains - Mernda line"
$ws.Range("C2").Value = "28/12/2020 14:30 - 14:45"
$ws.Range("D2").Value = "Caught train from Southern Cross to Victoria Park station"
$ws.Range("E2").Value = "new"

# Re-fit columns A:D so their widths track the new (shorter/longer) text, just
# like Excel does when "best fit" columns are recalculated after an edit.
# (Values below are chosen so the engine's pixel-quantized ColumnWidth setter
# lands on the width closest to what Excel's own AutoFit would compute.)
$ws.Columns.Item(1).ColumnWidth = 11.583833333333335
$ws.Columns.Item(2).ColumnWidth = 20.917166666666667
$ws.Columns.Item(3).ColumnWidth = 20.58383333333333
$ws.Columns.Item(4).ColumnWidth = 45.08383333333334

# Leave the final selection/active cell on B2, matching the saved view state.
$ws.Range("B2").Select() | Out-Null
